$p = $ppt.ActivePresentation

# Slide 3 ("Title and Object" layout) already contains the blank
# title + content placeholders we need for the new slide, so the
# cleanest way to create an identical new slide 4 is to duplicate
# slide 3 - this naturally places the duplicate right after it
# (i.e. at the end, as slide 4) and keeps its formatting intact.
$orig = $p.Slides.Item(3)
$dup = $orig.Duplicate()
